# "Generate Report for Handback"
#
# The localization-status report is regenerated once a handback has come in
# and is in sync with en-US:
#   - Overview/zh-cn/de-de "Status" cells flip from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - zh-cn + de-de sheets now have a "Latest Target File" (a.md, hyperlinked)
#     and a "Latest Handback File" filled in
#   - de-de additionally now has a "Latest Handback DateTime" (zh-cn's
#     handback datetime was already generated earlier)
#   - a couple of report columns get widened so the longer values fit

$wb = $excel.ActiveWorkbook

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e286523ba8130ea9104fd055de02b5186a4a308/e2e/a.md"

$hyperlinkColor = 15570276  # OLE (BGR) encoding of RGB FF6495ED, matching the workbook's HyperLink style

function Set-HyperlinkLook($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Overview sheet: Status text for both files, in both language columns
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the now-longer zh-cn / de-de status columns
$overview.Columns.Item(5).ColumnWidth = 29.083333333333336
$overview.Columns.Item(6).ColumnWidth = 29.083333333333336

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Target File
$zhcn.Range("I2").Value = "a.md"
$zhcn.Range("I3").Value = "a.md"

# Latest Handback File
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# Latest Handback DateTime (this language was already in sync before de-de)
$zhcn.Range("K2").Value = "2016-08-20 04:41:15"
$zhcn.Range("K3").Value = "2016-08-20 04:41:15"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $urlA, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $urlA, "", "", "a.md")
Set-HyperlinkLook $zhcn.Range("I2")
Set-HyperlinkLook $zhcn.Range("I3")

$zhcn.Columns.Item(3).ColumnWidth = 29.083333333333336
$zhcn.Columns.Item(10).ColumnWidth = 39.083333333333336

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Target File
$dede.Range("I2").Value = "a.md"
$dede.Range("I3").Value = "a.md"

# Latest Handback File
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# Latest Handback DateTime (just generated for de-de)
$dede.Range("K2").Value = "2016-08-20 04:41:21"
$dede.Range("K3").Value = "2016-08-20 04:41:21"

$dede.Hyperlinks.Add($dede.Range("I2"), $urlA, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $urlA, "", "", "a.md")
Set-HyperlinkLook $dede.Range("I2")
Set-HyperlinkLook $dede.Range("I3")

$dede.Columns.Item(3).ColumnWidth = 29.083333333333336
$dede.Columns.Item(10).ColumnWidth = 39.083333333333336
